# Auto-generated edit script: refreshes market-price columns (H-N)
# across multiple rows on sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
# per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 100000
$ws.Range("I18").Value = 100000
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 100000
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -99716
$ws.Range("H33").Value = 105.5
$ws.Range("I33").Value = 49.142857
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 49.142857
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = 179.857143
$ws.Range("H53").Value = 310.44446
$ws.Range("I53").Value = 98.5
$ws.Range("J53").Value = 371
$ws.Range("K53").Value = 98.5
$ws.Range("L53").Value = 371
$ws.Range("M53").Value = 538.5
$ws.Range("H99").Value = 125021150
$ws.Range("I99").Value = 27965.334
$ws.Range("J99").Value = 500000740
$ws.Range("K99").Value = 83896.00199999999
$ws.Range("L99").Value = 1500002220
$ws.Range("M99").Value = -82398.00199999999
$ws.Range("H132").Value = 4121.231
$ws.Range("I132").Value = 1984.174
$ws.Range("J132").Value = 20505.334
$ws.Range("K132").Value = 5952.522
$ws.Range("L132").Value = 61516.00199999999
$ws.Range("M132").Value = -3422.522
$ws.Range("H134").Value = 114607.2
$ws.Range("I134").Value = 179498
$ws.Range("J134").Value = 107397.11
$ws.Range("K134").Value = 179498
$ws.Range("L134").Value = 107397.11
$ws.Range("M134").Value = -174428
$ws.Range("N134").Value = -117537.11
$ws.Range("H137").Value = 2497.88
$ws.Range("I137").Value = 1356.9333
$ws.Range("J137").Value = 4209.3
$ws.Range("K137").Value = 4070.7999
$ws.Range("L137").Value = 12627.9
$ws.Range("M137").Value = -1520.7999
$ws.Range("N137").Value = -17727.9
$ws.Range("H138").Value = 2647.8
$ws.Range("I138").Value = 2085.8
$ws.Range("J138").Value = 3771.8
$ws.Range("K138").Value = 6257.400000000001
$ws.Range("L138").Value = 11315.4
$ws.Range("M138").Value = -1117.400000000001
$ws.Range("N138").Value = -21595.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2655.923
$ws.Range("I2").Value = 2113.8
$ws.Range("J2").Value = 4463
$ws.Range("K2").Value = 2113.8
$ws.Range("L2").Value = 4463
$ws.Range("M2").Value = -2000.8
$ws.Range("H32").Value = 5140.8296
$ws.Range("I32").Value = 4702.6665
$ws.Range("J32").Value = 14999.5
$ws.Range("K32").Value = 4702.6665
$ws.Range("L32").Value = 14999.5
$ws.Range("M32").Value = -4415.6665
$ws.Range("N32").Value = -15573.5
$ws.Range("H61").Value = 2180.74
$ws.Range("I61").Value = 2173.1428
$ws.Range("J61").Value = 2198.4666
$ws.Range("K61").Value = 2173.1428
$ws.Range("L61").Value = 2198.4666
$ws.Range("M61").Value = -1961.1428
$ws.Range("H116").Value = 2655.923
$ws.Range("I116").Value = 2113.8
$ws.Range("J116").Value = 4463
$ws.Range("K116").Value = 2113.8
$ws.Range("L116").Value = 4463
$ws.Range("M116").Value = 180.1999999999998
$ws.Range("H125").Value = 34999.4
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 34999.4
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 34999.4
$ws.Range("N125").Value = -44839.4
$ws.Range("H132").Value = 1540.0465
$ws.Range("I132").Value = 1541
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 4623
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -2093
$ws.Range("H136").Value = 2180.74
$ws.Range("I136").Value = 2173.1428
$ws.Range("J136").Value = 2198.4666
$ws.Range("K136").Value = 6519.428400000001
$ws.Range("L136").Value = 6595.399800000001
$ws.Range("M136").Value = -3969.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2655.923
$ws.Range("I3").Value = 2113.8
$ws.Range("J3").Value = 4463
$ws.Range("K3").Value = 2113.8
$ws.Range("L3").Value = 4463
$ws.Range("M3").Value = -1999.8
$ws.Range("H20").Value = 8230.370000000001
$ws.Range("I20").Value = 8609.182000000001
$ws.Range("J20").Value = 6563.6
$ws.Range("K20").Value = 8609.182000000001
$ws.Range("L20").Value = 6563.6
$ws.Range("M20").Value = -8362.182000000001
$ws.Range("N20").Value = -7057.6
$ws.Range("H94").Value = 1084.9048
$ws.Range("I94").Value = 853.2963
$ws.Range("J94").Value = 1501.8
$ws.Range("K94").Value = 853.2963
$ws.Range("L94").Value = 1501.8
$ws.Range("M94").Value = -402.2963
$ws.Range("H134").Value = 1528.45
$ws.Range("I134").Value = 1506.5625
$ws.Range("J134").Value = 1616
$ws.Range("K134").Value = 4519.6875
$ws.Range("L134").Value = 4848
$ws.Range("M134").Value = -1984.6875
$ws.Range("H140").Value = 134991.25
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 134991.25
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 134991.25
$ws.Range("N140").Value = -145351.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5959.1304
$ws.Range("I31").Value = 4011.5
$ws.Range("J31").Value = 6369.1577
$ws.Range("K31").Value = 4011.5
$ws.Range("L31").Value = 6369.1577
$ws.Range("M31").Value = -3716.5
$ws.Range("N31").Value = -6959.1577
$ws.Range("H34").Value = 5959.1304
$ws.Range("I34").Value = 4011.5
$ws.Range("J34").Value = 6369.1577
$ws.Range("K34").Value = 4011.5
$ws.Range("L34").Value = 6369.1577
$ws.Range("M34").Value = -3809.5
$ws.Range("N34").Value = -6773.1577
$ws.Range("H99").Value = 8449.76
$ws.Range("I99").Value = 11319.357
$ws.Range("J99").Value = 4797.5454
$ws.Range("K99").Value = 11319.357
$ws.Range("L99").Value = 4797.5454
$ws.Range("M99").Value = -9821.357
$ws.Range("N99").Value = -7793.5454
$ws.Range("H107").Value = 1767.7059
$ws.Range("I107").Value = 1445.1818
$ws.Range("J107").Value = 2359
$ws.Range("K107").Value = 1445.1818
$ws.Range("L107").Value = 2359
$ws.Range("M107").Value = 474.8181999999999
$ws.Range("H122").Value = 40015.777
$ws.Range("I122").Value = 63688.188
$ws.Range("J122").Value = 5583.1816
$ws.Range("K122").Value = 191064.564
$ws.Range("L122").Value = 16749.5448
$ws.Range("M122").Value = -188614.564
$ws.Range("H126").Value = 8449.76
$ws.Range("I126").Value = 11319.357
$ws.Range("J126").Value = 4797.5454
$ws.Range("K126").Value = 33958.071
$ws.Range("L126").Value = 14392.6362
$ws.Range("M126").Value = -31488.071
$ws.Range("N126").Value = -19332.6362
$ws.Range("H134").Value = 2371.5
$ws.Range("I134").Value = 1891.4474
$ws.Range("J134").Value = 5411.8335
$ws.Range("K134").Value = 5674.3422
$ws.Range("L134").Value = 16235.5005
$ws.Range("M134").Value = -3139.3422
$ws.Range("H141").Value = 177412.45
$ws.Range("I141").Value = 90000
$ws.Range("J141").Value = 182013.11
$ws.Range("K141").Value = 90000
$ws.Range("L141").Value = 182013.11
$ws.Range("M141").Value = -84820
$ws.Range("N141").Value = -192373.11

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2629.3845
$ws.Range("I75").Value = 665.3333
$ws.Range("J75").Value = 4312.857
$ws.Range("K75").Value = 1995.9999
$ws.Range("L75").Value = 12938.571
$ws.Range("M75").Value = -997.9999
$ws.Range("N75").Value = -14934.571
$ws.Range("H78").Value = 2629.3845
$ws.Range("I78").Value = 665.3333
$ws.Range("J78").Value = 4312.857
$ws.Range("K78").Value = 5987.9997
$ws.Range("L78").Value = 38815.713
$ws.Range("M78").Value = -995.9997000000003
$ws.Range("N78").Value = -48799.713

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 250296.16
$ws.Range("I2").Value = 500307.75
$ws.Range("J2").Value = 284.55
$ws.Range("K2").Value = 500307.75
$ws.Range("L2").Value = 284.55
$ws.Range("M2").Value = -500194.75
$ws.Range("N2").Value = -510.55
$ws.Range("H70").Value = 7836.25
$ws.Range("I70").Value = 8448.333000000001
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 8448.333000000001
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -8178.333000000001
$ws.Range("H73").Value = 7836.25
$ws.Range("I73").Value = 8448.333000000001
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 8448.333000000001
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -7512.333000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 9093
$ws.Range("I58").Value = 9093
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 9093
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -8833
$ws.Range("H136").Value = 3715.6667
$ws.Range("I136").Value = 2702.3809
$ws.Range("J136").Value = 6080
$ws.Range("K136").Value = 8107.1427
$ws.Range("L136").Value = 18240
$ws.Range("M136").Value = -5557.1427
$ws.Range("N136").Value = -23340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M63").ClearContents()
$ws.Range("H63").Value = 29999.666
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 29999.666
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 29999.666
$ws.Range("N63").Value = -31247.666
$ws.Range("M66").ClearContents()
$ws.Range("H66").Value = 29999.666
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 29999.666
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 89998.99800000001
$ws.Range("N66").Value = -96238.99800000001
$ws.Range("H70").Value = 18252.5
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 18252.5
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 18252.5
$ws.Range("N70").Value = -18882.5
$ws.Range("H73").Value = 18252.5
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 18252.5
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 18252.5
$ws.Range("N73").Value = -20436.5
$ws.Range("H122").Value = 1654.3429
$ws.Range("I122").Value = 1594.9231
$ws.Range("J122").Value = 1826
$ws.Range("K122").Value = 4784.7693
$ws.Range("L122").Value = 5478
$ws.Range("M122").Value = -2334.7693
$ws.Range("N122").Value = -10378
$ws.Range("H136").Value = 770.5333000000001
$ws.Range("I136").Value = 643.3461
$ws.Range("J136").Value = 1597.25
$ws.Range("K136").Value = 1930.0383
$ws.Range("L136").Value = 4791.75
$ws.Range("M136").Value = 619.9617000000001
